# Generate Report for Handback
#
# The 626b2942-f383-4f06-b035-e8b278fce277.md file has come back from
# localization in sync with en-US for both zh-cn and de-de. Update each
# language sheet's row for that file: flip its Status to "Handed back: in
# sync with en-US", stamp a Latest Handback DateTime, and record the
# Latest Target File / Latest Handback File hyperlinks for the file.

$wb = $excel.ActiveWorkbook

$langs = @(
    @{
        Sheet      = "zh-cn"
        Handback   = "2016-03-23 18:41:44"
        TargetUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/52e2014202e255bf55de3db600e9c54c04507d32/e2e/626b2942-f383-4f06-b035-e8b278fce277.md"
        HandbackUrl= "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/78b21aaaafe729a70fc15dbd9298f7da2baa3424/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/626b2942-f383-4f06-b035-e8b278fce277.412d32f0bcca43fd98a68bd24f35b5b9441fb961.zh-cn.xlf"
        XlfName    = "626b2942-f383-4f06-b035-e8b278fce277.412d32f0bcca43fd98a68bd24f35b5b9441fb961.zh-cn.xlf"
    },
    @{
        Sheet      = "de-de"
        Handback   = "2016-03-23 18:41:50"
        TargetUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/52e2014202e255bf55de3db600e9c54c04507d32/e2e/626b2942-f383-4f06-b035-e8b278fce277.md"
        HandbackUrl= "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/255dad796fd2e51a6af89c5c88d0ef95fed892f0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/626b2942-f383-4f06-b035-e8b278fce277.412d32f0bcca43fd98a68bd24f35b5b9441fb961.de-de.xlf"
        XlfName    = "626b2942-f383-4f06-b035-e8b278fce277.412d32f0bcca43fd98a68bd24f35b5b9441fb961.de-de.xlf"
    }
)

$mdName = "626b2942-f383-4f06-b035-e8b278fce277.md"

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Row 2 is the 626b2942-f383-4f06-b035-e8b278fce277.md file in every
    # per-language sheet. Status -> handed back, in sync with en-US.
    $ws.Range("C2").Value = "Handed back: in sync with en-US"

    # Latest Target File (F2) - the source .md that was targeted.
    $ws.Range("F2").Value = $mdName
    $ws.Hyperlinks.Add($ws.Range("F2"), $lang.TargetUrl, "", "", $mdName) | Out-Null

    # Latest Handback File (G2) - the localized .xlf handed back.
    $ws.Range("G2").Value = $lang.XlfName
    $ws.Hyperlinks.Add($ws.Range("G2"), $lang.HandbackUrl, "", "", $lang.XlfName) | Out-Null

    # Latest Handback DateTime (H2).
    $ws.Range("H2").Value = $lang.Handback
}
